$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four existing "press space bar to start session N" cells with
# the revised wording, and fill in the previously-missing "session 3" one.
# A leading apostrophe is used so the cells keep their existing
# quote-prefixed / centered-wrap style (matching their sibling cells)
# instead of Excel re-picking a plain style when the value is assigned.

# C12 (train1 column) - session 1
$ws.Range("C12").Value = "'세션 1`n`n'스페이스 바' 를 누르시면 세션 1 이 시작됩니다."

# E5 (test1 column) - session 3 (previously blank)
$ws.Range("E5").Value = "'세션 3`n`n'스페이스 바' 를 누르시면 세션 3 이 시작됩니다."

# F5 (test2 column) - session 4
$ws.Range("F5").Value = "'세션 4`n`n'스페이스 바' 를 누르시면 세션 4 가 시작됩니다."

# G9 (test3 column) - session 5
$ws.Range("G9").Value = "'세션 5`n`n'스페이스 바' 를 누르시면 세션 5 가 시작됩니다."

# D5 (train2 column) - session 2
$ws.Range("D5").Value = "'세션 2`n`n'스페이스 바' 를 누르시면 세션 2 가 시작됩니다."

# Move the view / selection from B8 to D5
$ws.Range("A4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D5").Select() | Out-Null
